$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'compression pants basketball men'
$ws.Range("A2").Value = 'basketball compression pants with pads'
$ws.Range("A3").Value = 'basketball black'
$ws.Range("A4").Value = 'sports compression pants men'
$ws.Range("A5").Value = 'knee protector men'
$ws.Range("A6").Value = 'athletic compression pants men'
$ws.Range("A7").Value = 'basketball tights for boys'
$ws.Range("A8").Value = 'protective knee pads for basketball'
$ws.Range("A9").Value = 'padded knee pads basketball youth'
$ws.Range("A10").Value = 'youth boys athletic leggings'
$ws.Range("A11").Value = 'baseball sliding shorts youth'
$ws.Range("A12").Value = 'gel knee pads for wrestling'
$ws.Range("A13").Value = 'running tights capri'
$ws.Range("A14").Value = 'bjj knee pads'
$ws.Range("A15").Value = 'football pants with pads adult'
$ws.Range("A16").Value = 'youth volleyball pads'
$ws.Range("A17").Value = 'mens compression tights black'
$ws.Range("A18").Value = 'little boy leggings for sports'
$ws.Range("A19").Value = 'baseball men pants'
$ws.Range("A20").Value = 'boys football pants with pads'
$ws.Range("A21").Value = 'black knee pads volleyball'
$ws.Range("A22").Value = 'knee compression for running'
$ws.Range("A23").Value = 'sweat pads men'
$ws.Range("A24").Value = 'basketball gear for men'
$ws.Range("A25").Value = 'baseball sliding shorts'
$ws.Range("A26").Value = 'black wrestling tights'
$ws.Range("A27").Value = 'mens basketballs'
$ws.Range("A28").Value = 'small black basketball'
$ws.Range("A29").Value = 'compression baseball'
$ws.Range("A30").Value = 'knee pads volleyball youth'
$ws.Range("A31").Value = 'baseball pants youth black'
$ws.Range("A32").Value = 'boys lacrosse pads'
$ws.Range("A33").Value = 'football tights for men'
$ws.Range("A34").Value = 'basketball youth'
$ws.Range("A35").Value = 'working knee pads'
$ws.Range("A36").Value = 'men compression clothes'
$ws.Range("A37").Value = 'male compression'
$ws.Range("A38").Value = 'kneepad volleyball'
$ws.Range("A39").Value = 'floor hockey pads'
$ws.Range("A40").Value = 'tactical pants men with knee pads'
$ws.Range("A41").Value = 'sports pants for men tall'
$ws.Range("A42").Value = 'lacrosse pads youth'
$ws.Range("A43").Value = 'mens capri shorts'
$ws.Range("A44").Value = 'sports leggings for boys'
$ws.Range("A45").Value = 'knee pads for washing floors'
$ws.Range("A46").Value = 'leggings under basketball shorts'
$ws.Range("A47").Value = 'basketball gym'
$ws.Range("A48").Value = 'capri pants mens'
$ws.Range("A49").Value = 'cycling tights for men'
$ws.Range("A50").Value = 'patella compression'
$ws.Range("A51").Value = 'hockey knee pad'
$ws.Range("A52").Value = 'mens gym pant'
$ws.Range("A53").Value = 'floor pads for gym'
$ws.Range("A54").Value = 'volleyball spandex'
$ws.Range("A55").Value = 'nike knee pads'
$ws.Range("A56").Value = 'mens volleyball knee pads'
$ws.Range("A57").Value = 'mcdavid knee pads'
$ws.Range("A58").Value = 'knee pad sleeve'
$ws.Range("A59").Value = 'cushy knee pads'
$ws.Range("A60").Value = 'men compression pants 3 4'
$ws.Range("A61").Value = 'men compression pants under armour'
$ws.Range("A62").Value = 'men compression pants nike'
$ws.Range("A63").Value = 'mens compression tights under armour'
$ws.Range("A64").Value = 'womens compression pants'
$ws.Range("A65").Value = 'men compression pants xxxl'
$ws.Range("A66").Value = 'men compression tights nike'
$ws.Range("A67").Value = 'mens compression tights short'
$ws.Range("A68").Value = 'mens compression tights thermal'
$ws.Range("A69").Value = 'mens compression pants white'
$ws.Range("A70").Value = 'mens compression pants navy blue'
$ws.Range("A71").Value = 'mens compression pants grey'
$ws.Range("A72").Value = 'mens compression pants nike'
$ws.Range("A73").Value = 'mens compression pants cold gear'
$ws.Range("A74").Value = 'eastbay mens compression tights'
$ws.Range("A75").Value = 'womens compression pant'
$ws.Range("A76").Value = 'ladies compression leggings'
$ws.Range("A77").Value = 'female compression pants'
$ws.Range("A78").Value = 'mcdavid basketball knee pads'
$ws.Range("A79").Value = 'robo knee pads'
$ws.Range("A80").Value = 'basketball long knee pads'
$ws.Range("A81").Value = 'basketball mcdavid knee pads'
$ws.Range("A82").Value = 'basketball nike knee pads'
$ws.Range("A83").Value = 'basketball tights for kids'
$ws.Range("A84").Value = 'ucla basketball youth'
$ws.Range("A85").Value = 'nba basketball youth'
$ws.Range("A86").Value = 'basketball 3 4 compression pants'
$ws.Range("A87").Value = 'basketball youth sleeve'
$ws.Range("A88").Value = 'basketball youth socks'
$ws.Range("A89").Value = 'knee pad nike'
$ws.Range("A90").Value = 'knee pad leggings'
$ws.Range("A91").Value = 'bicycle knee pads'
$ws.Range("A92").Value = 'rollerblade knee pads'
$ws.Range("A93").Value = 'black leggings pants'
$ws.Range("A94").Value = 'mtb knee pads'
$ws.Range("A95").Value = 'kids knee pads'
$ws.Range("A96").Value = 'knee pad set'
$ws.Range("A97").Value = 'compression shorts with knee pads'
$ws.Range("A98").Value = 'capri tights for women'
$ws.Range("A99").Value = 'knee pad biking'
$ws.Range("A100").Value = 'knee pads blue'
